# Fixed update to excel issue
#
# 1. Rename header "Requested quantity" -> "Weekly_PO_Qty" on "Weekly Quantity".
# 2. Rename header "Requested quantity" -> "Monthly_PO_Qty" on "Monthly Trend".
# 3. Add a new "PO Forecast" sheet (after "Monthly Trend") with forecast data.

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet after the last existing sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Headers
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match existing header style (bold / centered / bordered) used on the other sheets.
$wsWeekly.Range("A1:B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)

$data = @(
    @(44934.99999999999, 16, -0.09779519476231546, 30.5175760977157),
    @(44948.99999999999, 16, 1.936134972309997, 30.59432250146994),
    @(45011.99999999999, 17, 2.508564285648274, 32.90547549034741),
    @(45032.99999999999, 18, 3.947226980320172, 31.98357257510425),
    @(45039.99999999999, 18, 2.785475313520616, 31.95362048413374),
    @(45060.99999999999, 19, 4.34216257215018, 33.91082164608902),
    @(45067.99999999999, 19, 4.332095785972117, 32.98856411559374),
    @(45081.99999999999, 19, 5.226834162814677, 35.77133612546145),
    @(45088.99999999999, 19, 4.376124503698311, 33.73415381962095),
    @(45095.99999999999, 19, 3.89718018879696, 34.06976786651674),
    @(45102.99999999999, 19, 3.49402824134867, 34.816638154279),
    @(45109.99999999999, 20, 4.997032029896718, 34.36150918382533),
    @(45116.99999999999, 20, 5.652807813411259, 33.6543137993653),
    @(45130.99999999999, 20, 4.952292019739362, 34.13112665635121),
    @(45137.99999999999, 20, 5.283195219147673, 35.81748056881805),
    @(45151.99999999999, 21, 5.397858140124777, 34.83650719984336),
    @(45158.99999999999, 21, 6.079290489776307, 35.49536733304802),
    @(45165.99999999999, 21, 7.827125479791554, 36.03957524738408),
    @(45172.99999999999, 21, 6.46364533686632, 37.0147471059265),
    @(45179.99999999999, 21, 6.299777153682541, 36.18198519644645),
    @(45186.99999999999, 21, 7.199988969799558, 36.34894828931417),
    @(45193.99999999999, 21, 8.471067964058888, 36.21702766892296),
    @(45200.99999999999, 22, 6.273352772438047, 37.15463035124446),
    @(45207.99999999999, 22, 8.240718328949457, 36.70630594216842)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    for ($c = 0; $c -lt 4; $c++) {
        $wsForecast.Cells.Item($r + 2, $c + 1).Value = $data[$r][$c]
    }
}

# Match existing date-column style used for column A on the other sheets.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A25").PasteSpecial(-4122)

$wsForecast.Range("A1").Select()
